# CryCompanywiseStockReport_1.xlsx edit
#
# The report lists, per supplying company, stock-keeping units with
# Cost(D), MRP(E), Qty(F) and Value(G = D*F); each company block ends in a
# "Sub Total:" row whose B is SUM(G) over the block, and the sheet ends
# with an aggregate "Sub Total:" (= SUM of all block subtotals) followed
# by a "Grand Total:" that mirrors it. None of these are live formulas in
# the workbook - every cell is a stored literal - so quantity edits must
# be propagated by hand into G and every dependent subtotal/total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Quantity (column F) corrections -----------------------------------
# row number -> new quantity
$qtyEdits = @{
    11 = 16;
    25 = 34;
    32 = 80;
    81 = 55;
    89 = 69;
    98 = 246;
    106 = 57;
    114 = 325;
    139 = 42;
    141 = 33;
    146 = 155;
    170 = 12;
    173 = 9;
    178 = 10;
    219 = 6;
    221 = 38;
    247 = 200;
    257 = 9;
    285 = 23;
    288 = 26;
    295 = 40;
    296 = 103;
    297 = 131;
    313 = 166;
    325 = 56;
    329 = 52;
    330 = 128;
    338 = 14;
    341 = 233;
    391 = 46;
    394 = 41;
    400 = 142;
    402 = 71;
    411 = 65;
    412 = 69;
    413 = 49;
    416 = 80;
    417 = 74;
    419 = 67;
    427 = 2;
    440 = 350;
    445 = 144;
    451 = 309;
    455 = 61;
    499 = 650;
    502 = 144;
    503 = 151;
    516 = 42;
    524 = 64;
    551 = 27;
    568 = 17;
    570 = 28;
    578 = 91;
    579 = 37;
    606 = 180;
    621 = 6;
    669 = 7;
    671 = 80;
    672 = 16;
    687 = 68;
    690 = 13;
    692 = 83;
    697 = 59;
    699 = 124;
    701 = 111;
    747 = 1222;
    749 = 224;
    750 = 52;
    751 = 134
}

foreach ($row in $qtyEdits.Keys) {
    $qty  = $qtyEdits[$row]
    $cost = $ws.Cells.Item($row, 4).Value2      # column D
    $ws.Cells.Item($row, 6).Value = $qty                     # column F
    $ws.Cells.Item($row, 7).Value = [Math]::Round($cost * $qty, 10)   # column G
}

# --- 2. Two re-sequenced batch pairs ---------------------------------------
# Two items each have two stock-batch rows; the batches were re-ordered
# (later batch now listed first) and, for the HUL item, the moved batch's
# quantity was additionally trimmed by 3 units.
#   columns: B=code, D=cost, E=mrp, F=qty, G=value(=D*F)
function Set-Row($row, $code, $cost, $mrp, $qty) {
    $ws.Cells.Item($row, 2).Value = $code
    $ws.Cells.Item($row, 4).Value = $cost
    $ws.Cells.Item($row, 5).Value = $mrp
    $ws.Cells.Item($row, 6).Value = $qty
    $ws.Cells.Item($row, 7).Value = [Math]::Round($cost * $qty, 10)
}

# HUL-Bru Inst Poly 50g (rows 283/284): batch 61610 moves to row 283 with
# qty cut from 219 to 216; batch 57077 moves to row 284 unchanged.
Set-Row 283 61610 102.71 122.71 216
Set-Row 284 57077 93.08  111.2  1

# KUS-Floor Wiper (rows 396/397): batches simply swap rows, quantities
# unchanged.
Set-Row 396 58047 105.54 126.1 62
Set-Row 397 47097 112.28 134.16 15

# --- 3. Recompute every "Sub Total:" block and the grand total -------------
$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row

$blockStart = 1
$subtotalRows = New-Object System.Collections.Generic.List[int]
$blockTotal = 0.0

for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq "Sub Total:") {
        $ws.Cells.Item($r, 2).Value = [Math]::Round($blockTotal, 10)
        $subtotalRows.Add($r) | Out-Null
        $blockTotal = 0.0
        $blockStart = $r + 1
    } else {
        $g = $ws.Cells.Item($r, 7).Value2
        if ($g -is [double] -or $g -is [int]) {
            $blockTotal += $g
        }
    }
}

# The very last "Sub Total:" row aggregates all the per-company subtotals
# above it (not item rows), so redo it as the sum of the other subtotals.
$finalSubtotalRow = $subtotalRows[$subtotalRows.Count - 1]
$grandSum = 0.0
for ($i = 0; $i -lt $subtotalRows.Count - 1; $i++) {
    $grandSum += $ws.Cells.Item($subtotalRows[$i], 2).Value2
}
$ws.Cells.Item($finalSubtotalRow, 2).Value = [Math]::Round($grandSum, 10)

# "Grand Total:" row mirrors the final subtotal.
for ($r = $finalSubtotalRow + 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq "Grand Total:") {
        $ws.Cells.Item($r, 2).Value = [Math]::Round($grandSum, 10)
        break
    }
}
